$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") '25.965.99'
Set-TextValue $ws.Range("E2") '  +0.44%  '
Set-TextValue $ws.Range("D3") '1.587.70'
Set-TextValue $ws.Range("E3") '  +0.15%  '
Set-TextValue $ws.Range("E4") '  -0.17%  '
Set-TextValue $ws.Range("D5") '210.34'
Set-TextValue $ws.Range("E5") '  +0.16%  '
Set-TextValue $ws.Range("E6") '  -0.20%  '
Set-TextValue $ws.Range("E7") '  +0.09%  '
Set-TextValue $ws.Range("E8") '  -0.38%  '
Set-TextValue $ws.Range("E9") '  -1.04%  '
Set-TextValue $ws.Range("D10") '17.91'
Set-TextValue $ws.Range("E10") '  -0.95%  '
Set-TextValue $ws.Range("E11") '  +2.08%  '
Set-TextValue $ws.Range("D12") '1.809.50'
Set-TextValue $ws.Range("E12") '  +0.21%  '
Set-TextValue $ws.Range("D13") '1.588.39'
Set-TextValue $ws.Range("E13") '  +0.32%  '
Set-TextValue $ws.Range("E14") '  -1.25%  '
Set-TextValue $ws.Range("E15") '  +0.03%  '
Set-TextValue $ws.Range("D16") '25.954.52'
Set-TextValue $ws.Range("D17") '60.02'
Set-TextValue $ws.Range("E17") '  +0.43%  '
Set-TextValue $ws.Range("D18") '0.0₃0719'
Set-TextValue $ws.Range("E18") '  -0.43%  '
Set-TextValue $ws.Range("E19") '  -0.14%  '
Set-TextValue $ws.Range("D20") '199.35'
Set-TextValue $ws.Range("E20") '  +4.04%  '
Set-TextValue $ws.Range("E21") '  +0.62%  '
Set-TextValue $ws.Range("E22") '  -2.11%  '
Set-TextValue $ws.Range("E23") '  +0.45%  '
Set-TextValue $ws.Range("D24") '1.84'
Set-TextValue $ws.Range("E24") '  +8.34%  '
Set-TextValue $ws.Range("D25") '142.50'
Set-TextValue $ws.Range("E25") '  +0.31%  '
Set-TextValue $ws.Range("E26") '  -0.14%  '
Set-TextValue $ws.Range("E27") '  -8.51%  '
Set-TextValue $ws.Range("D28") '15.04'
Set-TextValue $ws.Range("E28") '  -0.38%  '
Set-TextValue $ws.Range("E29") '  -0.26%  '
Set-TextValue $ws.Range("E30") '  +0.11%  '
Set-TextValue $ws.Range("D31") '0.0473'
Set-TextValue $ws.Range("E31") '  +0.62%  '
Set-TextValue $ws.Range("E32") '  -0.03%  '
Set-TextValue $ws.Range("E33") '  -3.52%  '
Set-TextValue $ws.Range("E34") '  -1.84%  '
Set-TextValue $ws.Range("D35") '2.35'
Set-TextValue $ws.Range("E35") '  -0.21%  '
Set-TextValue $ws.Range("D36") '1.122.03'
Set-TextValue $ws.Range("E36") '  +1.84%  '
Set-TextValue $ws.Range("E37") '  +8.64%  '
Set-TextValue $ws.Range("E38") '  -0.17%  '
Set-TextValue $ws.Range("E39") '  -1.64%  '
Set-TextValue $ws.Range("E40") '  +0.58%  '
Set-TextValue $ws.Range("E41") '  -3.06%  '
Set-TextValue $ws.Range("D42") '0.777'
Set-TextValue $ws.Range("E42") '  -5.20%  '
Set-TextValue $ws.Range("D43") '1.720.93'
Set-TextValue $ws.Range("E43") '  +0.06%  '
Set-TextValue $ws.Range("E44") '  -1.57%  '
Set-TextValue $ws.Range("D45") '91.89'
Set-TextValue $ws.Range("E45") '  -2.09%  '
Set-TextValue $ws.Range("E46") '  -1.37%  '
Set-TextValue $ws.Range("D47") '53.21'
Set-TextValue $ws.Range("E47") '  -0.08%  '
Set-TextValue $ws.Range("E48") '  -1.18%  '
Set-TextValue $ws.Range("E49") '  -0.26%  '
Set-TextValue $ws.Range("E50") '  +0.07%  '
Set-TextValue $ws.Range("D51") '0.0₇0915'
Set-TextValue $ws.Range("E51") '  -14.64%  '
